# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# OFF sheet - Road ("R") row, updated Short/Deep Att/Comp + Int counts
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 159
$wsOff.Range("C3").Value = 111
$wsOff.Range("D3").Value = 42
$wsOff.Range("E3").Value = 22

# DEF sheet - Road ("R") row, updated Short/Deep Att/Comp + Int counts
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 198
$wsDef.Range("C3").Value = 142
$wsDef.Range("D3").Value = 65
$wsDef.Range("E3").Value = 35
